$d = $word.ActiveDocument

# 1. Merge the split runs in the "Adding 4th line for test." paragraph into
#    single runs: "Adding " + "4" -> "Adding 4", and " " + "line for test."
#    -> " line for test.". A Find/Replace of identical text over that span
#    normalizes/coalesces the adjacent plain runs as a side effect.
$d.Content.Find.Execute("Adding 4", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Adding 4", 2)

# 2. Insert a new blank paragraph, then a paragraph with the new text,
#    right after the "Adding 4th line for test." paragraph (i.e. before the
#    first of the pre-existing trailing blank paragraphs).
$trailingBlank = $d.Paragraphs.Item(5)
$trailingBlank.Range.InsertParagraphBefore()
$trailingBlank.Range.InsertParagraphBefore()

$d.Paragraphs.Item(6).Range.Text = "Testing git stash command"
